# Adição de mais 3 meses aos dados
# Adds "jun-jul-ago 2020", "jul-ago-set 2020" and "ago-set-out 2020" columns
# (AA, AB, AC) to the "Tabela" worksheet, extending the header row, the data
# row and the merged cell ranges that span the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Extend header row (row 4) with the three new quarter labels ---
$ws.Range("AA4").Value = "jun-jul-ago 2020"
$ws.Range("AB4").Value = "jul-ago-set 2020"
$ws.Range("AC4").Value = "ago-set-out 2020"

# --- Extend data row (row 5) with the three new values ---
$ws.Range("AA5").Value = 137510
$ws.Range("AB5").Value = 139961
$ws.Range("AC5").Value = 141908

# --- Re-span the merged header/footer ranges to cover the new columns ---
$ws.Range("A1:Z1").UnMerge()
$ws.Range("A1:AC1").Merge()

$ws.Range("A2:Z2").UnMerge()
$ws.Range("A2:AC2").Merge()

$ws.Range("B3:Z3").UnMerge()
$ws.Range("B3:AC3").Merge()

$ws.Range("A6:Z6").UnMerge()
$ws.Range("A6:AC6").Merge()
